$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "1.030", "0.07160") keep their exact literal digits/trailing zeros
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{ Row = 2; D = "27.493.57"; E = "  +2.13%  " }
    @{ Row = 3; D = "1.847.95"; E = "  +1.94%  " }
    @{ Row = 4; D = "1.030"; E = "  +2.76%  " }
    @{ Row = 5; D = "320.71"; E = "  +3.32%  " }
    @{ Row = 6; D = "1.025"; E = "  +2.23%  " }
    @{ Row = 7; D = "0.4376"; E = "  +2.12%  " }
    @{ Row = 8; D = "0.3768"; E = "  +2.18%  " }
    @{ Row = 9; D = "0.07393"; E = "  +2.24%  " }
    @{ Row = 10; D = "0.8730"; E = "  +1.35%  " }
    @{ Row = 11; D = "21.48"; E = "  +1.53%  " }
    @{ Row = 12; D = "1.858.11"; E = "  -7.77%  " }
    @{ Row = 13; D = "5.516"; E = "  +2.26%  " }
    @{ Row = 14; D = "6.686"; E = "  +0.82%  " }
    @{ Row = 15; D = "0.07160"; E = "  +3.92%  " }
    @{ Row = 16; D = "82.84"; E = "  +2.77%  " }
    @{ Row = 17; D = "1.032"; E = "  +2.87%  " }
    @{ Row = 18; D = "0.000009018"; E = "  +1.81%  " }
    @{ Row = 19; D = "1.026"; E = "  +2.20%  " }
    @{ Row = 20; D = "15.34"; E = "  +1.04%  " }
    @{ Row = 21; D = "27.512.41"; E = "  +2.05%  " }
    @{ Row = 22; D = "5.246"; E = "  +1.42%  " }
    @{ Row = 23; D = "11.31"; E = "  +2.51%  " }
    @{ Row = 24; D = "157.44"; E = "  +2.49%  " }
    @{ Row = 25; D = "1.908"; E = "  +1.42%  " }
    @{ Row = 26; D = "18.72"; E = "  +2.53%  " }
    @{ Row = 27; D = "1.969"; E = "  +5.08%  " }
    @{ Row = 28; D = "5.261"; E = "  +0.72%  " }
    @{ Row = 29; D = "117.03"; E = "  +1.95%  " }
    @{ Row = 30; D = "0.09037"; E = "  +1.04%  " }
    @{ Row = 31; D = "1.196"; E = "  +2.52%  " }
    @{ Row = 32; D = "0.7608"; E = "  +2.37%  " }
    @{ Row = 33; D = "4.516"; E = "  +2.18%  " }
    @{ Row = 34; D = "2.870"; E = "  +2.53%  " }
    @{ Row = 35; D = "1.027"; E = "  +1.88%  " }
    @{ Row = 36; D = "1.146"; E = "  +1.98%  " }
    @{ Row = 37; D = "0.01973"; E = "  +2.83%  " }
    @{ Row = 38; D = "0.05290"; E = "  +1.54%  " }
    @{ Row = 39; D = "0.5141"; E = "  +1.06%  " }
    @{ Row = 40; D = "2.792"; E = "  +3.03%  " }
    @{ Row = 41; D = "0.1672"; E = "  +1.91%  " }
    @{ Row = 42; D = "6.733"; E = "  +4.70%  " }
    @{ Row = 43; D = "8.472"; E = "  +2.34%  " }
    @{ Row = 44; D = "108.55"; E = "  +1.71%  " }
    @{ Row = 45; D = "10.58"; E = "  +1.72%  " }
    @{ Row = 46; D = "1.704"; E = "  +3.16%  " }
    @{ Row = 47; D = "0.4641"; E = "  +1.55%  " }
    @{ Row = 48; D = "0.06391"; E = "  +1.88%  " }
    @{ Row = 49; D = "1.846"; E = "  +2.17%  " }
    @{ Row = 50; D = "39.03"; E = "  +3.82%  " }
    @{ Row = 51; D = "63.84"; E = "  +0.36%  " }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
